$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.093281975671744
$ws.Range("C2").Value = 0.08507507989096297
$ws.Range("E2").Value = 0.08538992517543065
$ws.Range("F2").Value = 1.455836904019748
$ws.Range("G2").Value = 0.8581699763704052
$ws.Range("H2").Value = 0.01237382786354407
$ws.Range("I2").Value = 0.01829256405527779
$ws.Range("J2").Value = 0.6320295612235469
$ws.Range("K2").Value = 0.7063615233019291
$ws.Range("L2").Value = 0.0693021826315583
$ws.Range("M2").Value = 0.9332730324888985
$ws.Range("N2").Value = 0.3490029367378042

$ws.Range("B3").Value = 0.9524541764741912
$ws.Range("C3").Value = 0.07646643034472334
$ws.Range("E3").Value = 0.07648513219944775
$ws.Range("F3").Value = 1.368843988274563
$ws.Range("G3").Value = 0.8284924029238567
$ws.Range("H3").Value = 0.0159299741209209
$ws.Range("I3").Value = 0.02242137263576716
$ws.Range("J3").Value = 0.6230932810082521
$ws.Range("K3").Value = 0.6933520654533964
$ws.Range("L3").Value = 0.06590216712032948
$ws.Range("M3").Value = 0.8120117409995942
$ws.Range("N3").Value = 0.3069025293734171

$ws.Range("B4").Value = 0.8659112582748207
$ws.Range("C4").Value = 0.0712544047208965
$ws.Range("E4").Value = 0.07102412790389323
$ws.Range("F4").Value = 1.315852820557339
$ws.Range("G4").Value = 0.810575977605609
$ws.Range("H4").Value = 0.01842644085099576
$ws.Range("I4").Value = 0.02529843218986638
$ws.Range("J4").Value = 0.6178299002568508
$ws.Range("K4").Value = 0.6854657323808055
$ws.Range("L4").Value = 0.06378137117405558
$ws.Range("M4").Value = 0.7376966234987492
$ws.Range("N4").Value = 0.281212100777978

$ws.Range("B5").Value = 0.830401170257943
$ws.Range("C5").Value = 0.06926387922766253
$ws.Range("E5").Value = 0.06882364458397561
$ws.Range("F5").Value = 1.293066170664872
$ws.Range("G5").Value = 0.8020588224692347
$ws.Range("H5").Value = 0.01952633940361104
$ws.Range("I5").Value = 0.02665153249890384
$ws.Range("J5").Value = 0.6150722613960937
$ws.Range("K5").Value = 0.681316238068618
$ws.Range("L5").Value = 0.06294767684188152
$ws.Range("M5").Value = 0.7076663207514002
$ws.Range("N5").Value = 0.2711372611397564

$ws.Range("B6").Value = 0.8242363341401528
$ws.Range("C6").Value = 0.06907913160979007
$ws.Range("E6").Value = 0.06848778981706261
$ws.Range("F6").Value = 1.287710984612417
$ws.Range("G6").Value = 0.7990770952719686
$ws.Range("H6").Value = 0.01972164168856838
$ws.Range("I6").Value = 0.02700777306797875
$ws.Range("J6").Value = 0.61380360242498
$ws.Range("K6").Value = 0.6794576290659187
$ws.Range("L6").Value = 0.06285740086205038
$ws.Range("M6").Value = 0.7029598731857902
$ws.Range("N6").Value = 0.2699051309896134

$ws.Range("B7").Value = 0.8647001996645827
$ws.Range("C7").Value = 0.07162404722142668
$ws.Range("E7").Value = 0.07107511908766995
$ws.Range("F7").Value = 1.311237758358118
$ws.Range("G7").Value = 0.8061698857796813
$ws.Range("H7").Value = 0.01846207730277794
$ws.Range("I7").Value = 0.02564383799536252
$ws.Range("J7").Value = 0.6155728921965533
$ws.Range("K7").Value = 0.682213939165166
$ws.Range("L7").Value = 0.06390406420799444
$ws.Range("M7").Value = 0.7380537046250026
$ws.Range("N7").Value = 0.2822762367261191

$ws.Range("B8").Value = 1.043750003761232
$ws.Range("C8").Value = 0.08262005001881079
$ws.Range("E8").Value = 0.08242482740810075
$ws.Range("F8").Value = 1.420053827875407
$ws.Range("G8").Value = 0.8422096223746678
$ws.Range("H8").Value = 0.01355598165122232
$ws.Range("I8").Value = 0.0200338621903926
$ws.Range("J8").Value = 0.6259831509532034
$ws.Range("K8").Value = 0.6976455919231768
$ws.Range("L8").Value = 0.06831968788721987
$ws.Range("M8").Value = 0.8924399008360524
$ws.Range("N8").Value = 0.3360500479423507

$ws.Range("B9").Value = 1.395959224743308
$ws.Range("C9").Value = 0.1041923927270147
$ws.Range("E9").Value = 0.1046284052542426
$ws.Range("F9").Value = 1.64492738871968
$ws.Range("G9").Value = 0.9232644508667534
$ws.Range("H9").Value = 0.006522965406169035
$ws.Range("I9").Value = 0.01139446506336572
$ws.Range("J9").Value = 0.6522366951101617
$ws.Range("K9").Value = 0.7345043330898733
$ws.Range("L9").Value = 0.07648423476448585
$ws.Range("M9").Value = 1.195375516499809
$ws.Range("N9").Value = 0.4407837205027363

$ws.Range("B10").Value = 1.652680282774242
$ws.Range("C10").Value = 0.1219032078129985
$ws.Range("E10").Value = 0.117431541101169
$ws.Range("F10").Value = 1.781538901661818
$ws.Range("G10").Value = 0.9650617197813318
$ws.Range("H10").Value = 0.003477248552416601
$ws.Range("I10").Value = 0.007227714249262718
$ws.Range("J10").Value = 0.6627681373475269
$ws.Range("K10").Value = 0.7482792870827524
$ws.Range("L10").Value = 0.08367710291179264
$ws.Range("M10").Value = 1.421698833063829
$ws.Range("N10").Value = 0.5061994491579185

$ws.Range("B11").Value = 1.750283911314966
$ws.Range("C11").Value = 0.1404128691112305
$ws.Range("E11").Value = 0.09380719289734074
$ws.Range("F11").Value = 1.589788940347134
$ws.Range("G11").Value = 0.830649117284338
$ws.Range("H11").Value = 0.0219262222511496
$ws.Range("I11").Value = 0.006833127028532893
$ws.Range("J11").Value = 0.5894533028256177
$ws.Range("K11").Value = 0.6476062182469704
$ws.Range("L11").Value = 0.1027467878698118
$ws.Range("M11").Value = 1.542140204595569
$ws.Range("N11").Value = 0.427452925197116

$ws.Range("B12").Value = 1.779136072704944
$ws.Range("C12").Value = 0.1526045510101568
$ws.Range("E12").Value = 0.07403014774544658
$ws.Range("F12").Value = 1.414571805655399
$ws.Range("G12").Value = 0.722768263584058
$ws.Range("H12").Value = 0.06079800761126819
$ws.Range("I12").Value = 0.006698870562788173
$ws.Range("J12").Value = 0.5323827291311858
$ws.Range("K12").Value = 0.5714355439599075
$ws.Range("L12").Value = 0.1231656216671126
$ws.Range("M12").Value = 1.595330444219741
$ws.Range("N12").Value = 0.3517995510593863

$ws.Range("B13").Value = 1.755253020276285
$ws.Range("C13").Value = 0.1610958403534113
$ws.Range("E13").Value = 0.05628254943498945
$ws.Range("F13").Value = 1.23565916559879
$ws.Range("G13").Value = 0.62226777618271
$ws.Range("H13").Value = 0.1170568942675771
$ws.Range("I13").Value = 0.0071624599489164
$ws.Range("J13").Value = 0.4809433539438004
$ws.Range("K13").Value = 0.5046371066483353
$ws.Range("L13").Value = 0.1457164404927838
$ws.Range("M13").Value = 1.600464229348802
$ws.Range("N13").Value = 0.2765443478965892

$ws.Range("B14").Value = 1.71448697441835
$ws.Range("C14").Value = 0.1655899676604093
$ws.Range("E14").Value = 0.04530126488738517
$ws.Range("F14").Value = 1.110004711076868
$ws.Range("G14").Value = 0.5558395884383316
$ws.Range("H14").Value = 0.1669100932534775
$ws.Range("I14").Value = 0.007878577751207594
$ws.Range("J14").Value = 0.4476940684529751
$ws.Range("K14").Value = 0.4626333175512833
$ws.Range("L14").Value = 0.1631918027038068
$ws.Range("M14").Value = 1.582532637840956
$ws.Range("N14").Value = 0.2257154683228748

$ws.Range("B15").Value = 1.693097561182555
$ws.Range("C15").Value = 0.1659813807862491
$ws.Range("E15").Value = 0.04271233018342713
$ws.Range("F15").Value = 1.076500683241903
$ws.Range("G15").Value = 0.5392208407508576
$ws.Range("H15").Value = 0.1796611713189691
$ws.Range("I15").Value = 0.008339240464145803
$ws.Range("J15").Value = 0.4397853445093318
$ws.Range("K15").Value = 0.4528219008124559
$ws.Range("L15").Value = 0.1672723193457983
$ws.Range("M15").Value = 1.568129793861317
$ws.Range("N15").Value = 0.2130076422416494

$ws.Range("B16").Value = 1.587092312751537
$ws.Range("C16").Value = 0.1564589594680967
$ws.Range("E16").Value = 0.04136968877471681
$ws.Range("F16").Value = 1.054939060740637
$ws.Range("G16").Value = 0.5402063344071522
$ws.Range("H16").Value = 0.1682136966891079
$ws.Range("I16").Value = 0.01010386163042032
$ws.Range("J16").Value = 0.4446513501943201
$ws.Range("K16").Value = 0.459252140340471
$ws.Range("L16").Value = 0.1589881654624676
$ws.Range("M16").Value = 1.469179427425871
$ws.Range("N16").Value = 0.2030179333296047

$ws.Range("B17").Value = 1.527576255672585
$ws.Range("C17").Value = 0.1467907338459042
$ws.Range("E17").Value = 0.04581719079581603
$ws.Range("F17").Value = 1.103664078110569
$ws.Range("G17").Value = 0.5752788813879448
$ws.Range("H17").Value = 0.1311839259339251
$ws.Range("I17").Value = 0.01113440001121813
$ws.Range("J17").Value = 0.4656453396998188
$ws.Range("K17").Value = 0.4857119112739099
$ws.Range("L17").Value = 0.1429720921436726
$ws.Range("M17").Value = 1.403107930682438
$ws.Range("N17").Value = 0.2214387166748821

$ws.Range("B18").Value = 1.501548410291008
$ws.Range("C18").Value = 0.1358147636987184
$ws.Range("E18").Value = 0.05725351574883852
$ws.Range("F18").Value = 1.226314707814922
$ws.Range("G18").Value = 0.6504752523075297
$ws.Range("H18").Value = 0.07849328706546288
$ws.Range("I18").Value = 0.0111793897325656
$ws.Range("J18").Value = 0.5062526156416425
$ws.Range("K18").Value = 0.5378516535648821
$ws.Range("L18").Value = 0.1208840986101691
$ws.Range("M18").Value = 1.357097132244149
$ws.Range("N18").Value = 0.2692208115304311

$ws.Range("B19").Value = 1.502652139109614
$ws.Range("C19").Value = 0.1261946551345545
$ws.Range("E19").Value = 0.07597771784481822
$ws.Range("F19").Value = 1.401735362301835
$ws.Range("G19").Value = 0.7533140669749088
$ws.Range("H19").Value = 0.03266087687700292
$ws.Range("I19").Value = 0.01098055394270236
$ws.Range("J19").Value = 0.5597972853401387
$ws.Range("K19").Value = 0.6083314591073901
$ws.Range("L19").Value = 0.1004138900527529
$ws.Range("M19").Value = 1.332262674065021
$ws.Range("N19").Value = 0.3446594484499599

$ws.Range("B20").Value = 1.582737368223263
$ws.Range("C20").Value = 0.1185422324193723
$ws.Range("E20").Value = 0.1141603513264435
$ws.Range("F20").Value = 1.730937468130634
$ws.Range("G20").Value = 0.9400887606990267
$ws.Range("H20").Value = 0.004172858083301367
$ws.Range("I20").Value = 0.009107693756071455
$ws.Range("J20").Value = 0.6528184487745818
$ws.Range("K20").Value = 0.7344068369858903
$ws.Range("L20").Value = 0.08234114241101054
$ws.Range("M20").Value = 1.364614790954477
$ws.Range("N20").Value = 0.4921997698551763

$ws.Range("B21").Value = 1.783633601211335
$ws.Range("C21").Value = 0.1305767459955618
$ws.Range("E21").Value = 0.1296525094859113
$ws.Range("F21").Value = 1.885211187958063
$ws.Range("G21").Value = 1.003280006304152
$ws.Range("H21").Value = 0.001948131423379795
$ws.Range("I21").Value = 0.006240099897292239
$ws.Range("J21").Value = 0.6768211774315063
$ws.Range("K21").Value = 0.7666313960427686
$ws.Range("L21").Value = 0.08579353012170898
$ws.Range("M21").Value = 1.535755839401133
$ws.Range("N21").Value = 0.5626058171241226

$ws.Range("B22").Value = 1.915163484369202
$ws.Range("C22").Value = 0.138607658392857
$ws.Range("E22").Value = 0.1378456465574445
$ws.Range("F22").Value = 1.977411176289607
$ws.Range("G22").Value = 1.041572242001166
$ws.Range("H22").Value = 0.001036090207978013
$ws.Range("I22").Value = 0.004479066344047489
$ws.Range("J22").Value = 0.6910223450461217
$ws.Range("K22").Value = 0.7855984329363181
$ws.Range("L22").Value = 0.08848152967723877
$ws.Range("M22").Value = 1.648028129793005
$ws.Range("N22").Value = 0.600538682179689

$ws.Range("B23").Value = 1.845860702267998
$ws.Range("C23").Value = 0.1338103565163706
$ws.Range("E23").Value = 0.1333725032596185
$ws.Range("F23").Value = 1.933182054139195
$ws.Range("G23").Value = 1.026049922576391
$ws.Range("H23").Value = 0.001479980289839755
$ws.Range("I23").Value = 0.005037574488074625
$ws.Range("J23").Value = 0.6859498113528275
$ws.Range("K23").Value = 0.779127824767663
$ws.Range("L23").Value = 0.08686599214568602
$ws.Range("M23").Value = 1.587145936949156
$ws.Range("N23").Value = 0.5788073445441171

$ws.Range("B24").Value = 1.582469493656816
$ws.Range("C24").Value = 0.1167902925550166
$ws.Range("E24").Value = 0.1166628734620119
$ws.Range("F24").Value = 1.759406336003025
$ws.Range("G24").Value = 0.9607057200081641
$ws.Range("H24").Value = 0.003951972203180776
$ws.Range("I24").Value = 0.008523363479961077
$ws.Range("J24").Value = 0.663464136643185
$ws.Range("K24").Value = 0.7493472573247146
$ws.Range("L24").Value = 0.08100165223104128
$ws.Range("M24").Value = 1.359035914590208
$ws.Range("N24").Value = 0.4993656055054601

$ws.Range("B25").Value = 1.299275490566828
$ws.Range("C25").Value = 0.09899594044113513
$ws.Range("E25").Value = 0.09875305745961072
$ws.Range("F25").Value = 1.575858076795782
$ws.Range("G25").Value = 0.8933164794705419
$ws.Range("H25").Value = 0.008148279824538196
$ws.Range("I25").Value = 0.01393217447012951
$ws.Range("J25").Value = 0.6409683330169713
$ws.Range("K25").Value = 0.7187232638258934
$ws.Range("L25").Value = 0.07456280991243247
$ws.Range("M25").Value = 1.114578506080477
$ws.Range("N25").Value = 0.4144741043151328
